$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so the updated
# values are written back as plain strings (matching original formatting),
# not auto-converted to numbers/percentages.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "327.14"
$ws.Range("E2").Value = "-1.08%"
$ws.Range("D3").Value = "43.74"
$ws.Range("E3").Value = "5.26%"
$ws.Range("D4").Value = "5.464"
$ws.Range("E4").Value = "-3.93%"
$ws.Range("D5").Value = "0.08070"
$ws.Range("E5").Value = "-4.19%"
$ws.Range("D6").Value = "8.651"
$ws.Range("E6").Value = "-2.00%"
$ws.Range("D7").Value = "4.294"
$ws.Range("E7").Value = "-4.18%"
$ws.Range("D8").Value = "1.882"
$ws.Range("E8").Value = "-6.00%"
$ws.Range("E9").Value = "-8.35%"
$ws.Range("D10").Value = "0.9369"
$ws.Range("E10").Value = "1.14%"
$ws.Range("D11").Value = "0.1191"
$ws.Range("E11").Value = "-6.73%"
$ws.Range("D12").Value = "0.1897"
$ws.Range("E12").Value = "-3.43%"
$ws.Range("D13").Value = "0.09484"
$ws.Range("E13").Value = "1.13%"
$ws.Range("E14").Value = "2.10%"
$ws.Range("E15").Value = "0.46%"
$ws.Range("D16").Value = "0.001294"
$ws.Range("E16").Value = "0.16%"
$ws.Range("D17").Value = "0.005962"
$ws.Range("E17").Value = "-2.52%"
$ws.Range("E18").Value = "4.44%"
$ws.Range("E19").Value = "-0.04%"
$ws.Range("D20").Value = "8.624"
$ws.Range("E20").Value = "-3.68%"
$ws.Range("D21").Value = "0.1366"
$ws.Range("E21").Value = "0.21%"
$ws.Range("E22").Value = "-0.64%"
$ws.Range("D23").Value = "0.04362"
$ws.Range("E23").Value = "-1.29%"
$ws.Range("E24").Value = "-0.82%"
$ws.Range("D25").Value = "0.004342"
$ws.Range("E26").Value = "3.44%"
$ws.Range("D27").Value = "0.0004003"
$ws.Range("E27").Value = "0.13%"
$ws.Range("D39").Value = "0.02666"
$ws.Range("E39").Value = "-5.59%"
$ws.Range("D40").Value = "0.05417"
$ws.Range("E40").Value = "-1.92%"
$ws.Range("D41").Value = "0.007633"
$ws.Range("E41").Value = "-3.65%"
$ws.Range("D42").Value = "0.01045"
$ws.Range("E42").Value = "16.44%"
$ws.Range("D43").Value = "0.1389"
$ws.Range("E43").Value = "-3.47%"
$ws.Range("D44").Value = "0.002099"
$ws.Range("E44").Value = "0.20%"
$ws.Range("D45").Value = "0.009707"
$ws.Range("E45").Value = "-15.19%"
$ws.Range("D46").Value = "0.00006878"
$ws.Range("E46").Value = "-1.21%"
$ws.Range("E47").Value = "0.09%"
$ws.Range("D48").Value = "0.003560"
$ws.Range("E48").Value = "9.03%"
$ws.Range("D49").Value = "0.002278"
$ws.Range("E49").Value = "-0.20%"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").Value = "0.09%"
$ws.Range("E51").Value = "0.09%"

# Remove the temporary text-number-format so the cell styling matches
# the original (unstyled) cells.
$ws.Range("D2:E51").ClearFormats()
